$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1: thin top+bottom border only (no left/right), default (non-bold) style.
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none
$c1.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none

# D1: thin top+right+bottom border (no left), default (non-bold) style.
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Reuse the exact formatting built above (copy/paste formats only) so the
# workbook ends up with the same two new cell styles shared across both
# sheets, instead of synthesizing duplicate style/border definitions.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
